# Apply updated crypto price/volume figures per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values must stay plain text exactly as scraped -
# force Text format before the write so Excel does not auto-coerce
# numeric-looking strings (e.g. "230.76") into actual numbers, then
# restore the default style so no stray number-format is left behind.
function Set-TextValue($rangeRef, $text) {
    $rng = $ws.Range($rangeRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "43.644.91"
$ws.Range("E2").Value = "  +4.13%  "
Set-TextValue "D3" "2.263.54"
$ws.Range("E3").Value = "  +1.34%  "
$ws.Range("E4").Value = "  -0.14%  "
Set-TextValue "D5" "230.76"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("E6").Value = "  -0.02%  "
Set-TextValue "D7" "61.34"
$ws.Range("E7").Value = "  -0.71%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  +4.56%  "
Set-TextValue "D10" "58.07"
$ws.Range("E10").Value = "  -2.17%  "
Set-TextValue "D11" "0.0938"
$ws.Range("E11").Value = "  +5.20%  "
$ws.Range("E12").Value = "  +0.60%  "
Set-TextValue "D13" "2.598.84"
$ws.Range("E13").Value = "  +1.25%  "
Set-TextValue "D14" "15.61"
$ws.Range("E14").Value = "  -0.32%  "
Set-TextValue "D15" "23.64"
$ws.Range("E15").Value = "  +7.39%  "
Set-TextValue "D16" "5.79"
$ws.Range("E16").Value = "  +3.62%  "
Set-TextValue "D17" "0.810"
$ws.Range("E17").Value = "  +0.93%  "
Set-TextValue "D18" "2.264.46"
$ws.Range("E18").Value = "  +1.16%  "
Set-TextValue "D19" "43.541.06"
$ws.Range("E19").Value = "  +4.25%  "
Set-TextValue "D20" "0.0₃0934"
$ws.Range("E20").Value = "  +4.10%  "
Set-TextValue "D21" "72.86"
$ws.Range("E21").Value = "  +1.03%  "
Set-TextValue "D22" "6.22"
$ws.Range("E22").Value = "  +2.74%  "
Set-TextValue "D23" "252.65"
$ws.Range("E23").Value = "  +1.07%  "
Set-TextValue "D25" "2.55"
$ws.Range("E25").Value = "  +6.52%  "
$ws.Range("E26").Value = "  -0.22%  "
Set-TextValue "D27" "9.84"
$ws.Range("E27").Value = "  +1.53%  "
Set-TextValue "D28" "170.84"
$ws.Range("E28").Value = "  +2.51%  "
$ws.Range("E29").Value = "  -1.57%  "
Set-TextValue "D30" "20.49"
$ws.Range("E30").Value = "  +2.61%  "
Set-TextValue "D31" "1.44"
$ws.Range("E31").Value = "  +1.97%  "
$ws.Range("E32").Value = "  +0.75%  "
$ws.Range("E33").Value = "  -0.03%  "
Set-TextValue "D34" "5.07"
$ws.Range("E34").Value = "  +1.11%  "
Set-TextValue "D35" "4.80"
$ws.Range("E35").Value = "  +2.25%  "
Set-TextValue "D36" "0.0659"
$ws.Range("E36").Value = "  +3.74%  "
Set-TextValue "D37" "6.45"
$ws.Range("E37").Value = "  -3.08%  "
Set-TextValue "D38" "2.40"
$ws.Range("E38").Value = "  +1.10%  "
Set-TextValue "D39" "3.59"
$ws.Range("E39").Value = "  -1.88%  "
$ws.Range("E40").Value = "  +3.96%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  -10.44%  "
Set-TextValue "D43" "8.70"
$ws.Range("E43").Value = "  +1.43%  "
Set-TextValue "D44" "0.0995"
$ws.Range("E44").Value = "  +1.84%  "
$ws.Range("E45").Value = "  -6.76%  "
$ws.Range("E46").Value = "  -0.84%  "
Set-TextValue "D47" "98.08"
$ws.Range("E47").Value = "  -0.81%  "
Set-TextValue "D48" "1.471.83"
$ws.Range("E48").Value = "  -0.51%  "
Set-TextValue "D49" "16.62"
$ws.Range("E49").Value = "  +0.49%  "
$ws.Range("E50").Value = "  +0.52%  "
Set-TextValue "D51" "2.27"
$ws.Range("E51").Value = "  +8.01%  "
